# Updated symbol list on Mon Dec 19 15:43:21 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.41"
$ws.Range("D3").Value = "'21.71"
$ws.Range("D4").Value = "'5.425"
$ws.Range("D5").Value = "'0.05693"
$ws.Range("D7").Value = "'0.8092"
$ws.Range("D8").Value = "'1.023"
$ws.Range("D10").Value = "'0.07535"
$ws.Range("D12").Value = "'0.03052"
$ws.Range("D13").Value = "'0.09267"
$ws.Range("D14").Value = "'3.616"
$ws.Range("E14").Value = "13MCDexMCBBestin24h"
$ws.Range("D15").Value = "'0.001657"
$ws.Range("D16").Value = "'0.04711"
$ws.Range("D17").Value = "'0.0005852"
$ws.Range("E17").Value = "16OneONEWorstin24h"
$ws.Range("D18").Value = "'0.006349"
$ws.Range("D19").Value = "'0.005033"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("D20").Value = "'0.001043"
$ws.Range("D21").Value = "'0.0001501"
$ws.Range("D23").Value = "'3.769"
$ws.Range("D24").Value = "'6.411"
$ws.Range("D26").Value = "'0.3289"
$ws.Range("D27").Value = "'0.1335"
$ws.Range("D40").Value = "'0.04054"
$ws.Range("D41").Value = "'0.006969"
$ws.Range("D42").Value = "'0.1042"
$ws.Range("D43").Value = "'0.003503"
$ws.Range("D44").Value = "'0.008522"
$ws.Range("D45").Value = "'0.00005942"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D47").Value = "'0.0005502"
$ws.Range("E47").Value = "46ACDXExchangeACXT"
$ws.Range("D49").Value = "'0.007718"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("D51").Value = "'0.01011"
